$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "water loss"
$ws.Range("A15").Value = "drinking water treatment"

$ws.Range("H18").Select()
